$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102
$ws.Range("B102").Value = 7494646
$ws.Range("E102").Value = 'OHiggins'
$ws.Range("F102").Value = 'Cobresal'
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 3
$ws.Range("N102").Value = 2.3
$ws.Range("O102").Value = 2.1
$ws.Range("P102").Value = 3.5
$ws.Range("Q102").Value = 3.5
$ws.Range("R102").Value = -0.25
$ws.Range("S102").Value = 1.8
$ws.Range("T102").Value = 2.05
$ws.Range("Y102").Value = 2.5
$ws.Range("AA102").Value = -0.5
$ws.Range("AB102").Value = 0.5249999999999999

# Row 103
$ws.Range("B103").Value = 7494647
$ws.Range("E103").Value = 'Huachipato'
$ws.Range("F103").Value = 'Universidad Catolica'
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 1
$ws.Range("J103").Value = 1
$ws.Range("L103").Value = 2.2
$ws.Range("N103").Value = 3.2
$ws.Range("O103").Value = 1.8
$ws.Range("P103").Value = 3.6
$ws.Range("Q103").Value = 4.333
$ws.Range("R103").Value = -0.75
$ws.Range("S103").Value = 1.975
$ws.Range("T103").Value = 1.875
$ws.Range("Y103").Value = 2.6
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.875

# Row 115
$ws.Range("B115").Value = 6078997
$ws.Range("E115").Value = 'Union Espanola'
$ws.Range("F115").Value = 'Cobresal'
$ws.Range("G115").Value = 1
$ws.Range("L115").Value = 3.8
$ws.Range("M115").Value = 3.6
$ws.Range("N115").Value = 1.909
$ws.Range("O115").Value = 2.7
$ws.Range("P115").Value = 3.6
$ws.Range("Q115").Value = 2.45
$ws.Range("R115").Value = 0
$ws.Range("S115").Value = 1.975
$ws.Range("V115").Value = 1.775
$ws.Range("W115").Value = 2.025
$ws.Range("X115").Value = 1.7
$ws.Range("AA115").Value = 0.9750000000000001
$ws.Range("AD115").Value = 1.025

# Row 116
$ws.Range("B116").Value = 6143704
$ws.Range("E116").Value = 'Curico Unido'
$ws.Range("F116").Value = 'Colo Colo'
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 1
$ws.Range("K116").Value = 'A'
$ws.Range("L116").Value = 6.5
$ws.Range("M116").Value = 4.75
$ws.Range("N116").Value = 1.4
$ws.Range("O116").Value = 12
$ws.Range("P116").Value = 8.5
$ws.Range("Q116").Value = 1.166
$ws.Range("R116").Value = 2
$ws.Range("S116").Value = 2
$ws.Range("T116").Value = 1.8
$ws.Range("U116").Value = 3.25
$ws.Range("V116").Value = 1.875
$ws.Range("W116").Value = 1.925
$ws.Range("X116").Value = -1
$ws.Range("Z116").Value = 0.1659999999999999
$ws.Range("AA116").Value = 1
$ws.Range("AD116").Value = 0.925

# Row 117
$ws.Range("B117").Value = 6078267
$ws.Range("E117").Value = 'Huachipato'
$ws.Range("F117").Value = 'Audax Italiano'
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 0
$ws.Range("K117").Value = 'H'
$ws.Range("L117").Value = 1.5
$ws.Range("M117").Value = 4.333
$ws.Range("N117").Value = 6
$ws.Range("O117").Value = 1.444
$ws.Range("P117").Value = 4.75
$ws.Range("Q117").Value = 7
$ws.Range("R117").Value = -1.25
$ws.Range("S117").Value = 2.025
$ws.Range("T117").Value = 1.825
$ws.Range("U117").Value = 2.75
$ws.Range("V117").Value = 1.8
$ws.Range("W117").Value = 2.05
$ws.Range("X117").Value = 0.444
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 1.025
$ws.Range("AD117").Value = 1.05

# Row 137
$ws.Range("B137").Value = 7723533
$ws.Range("E137").Value = 'OHiggins'
$ws.Range("F137").Value = 'Everton de Vina'
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 1
$ws.Range("K137").Value = 'H'
$ws.Range("L137").Value = 3
$ws.Range("M137").Value = 3.2
$ws.Range("N137").Value = 2.375
$ws.Range("O137").Value = 2.3
$ws.Range("P137").Value = 3.1
$ws.Range("Q137").Value = 3.3
$ws.Range("R137").Value = -0.25
$ws.Range("S137").Value = 1.9
$ws.Range("T137").Value = 1.9
$ws.Range("V137").Value = 1.95
$ws.Range("W137").Value = 1.85
$ws.Range("X137").Value = 1.3
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = 0.8999999999999999
$ws.Range("AB137").Value = -1
$ws.Range("AC137").Value = 0.95
$ws.Range("AD137").Value = -1

# Row 138
$ws.Range("B138").Value = 7723528
$ws.Range("E138").Value = 'Palestino'
$ws.Range("F138").Value = 'Universidad Catolica'
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 'A'
$ws.Range("L138").Value = 1.95
$ws.Range("M138").Value = 3.6
$ws.Range("N138").Value = 3.4
$ws.Range("O138").Value = 2.375
$ws.Range("P138").Value = 3.5
$ws.Range("Q138").Value = 2.9
$ws.Range("R138").Value = 0
$ws.Range("S138").Value = 1.8
$ws.Range("T138").Value = 2.05
$ws.Range("V138").Value = 1.8
$ws.Range("W138").Value = 2.05
$ws.Range("X138").Value = -1
$ws.Range("Z138").Value = 1.9
$ws.Range("AA138").Value = -1
$ws.Range("AB138").Value = 1.05
$ws.Range("AC138").Value = -0.5
$ws.Range("AD138").Value = 0.5249999999999999
